$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 ("Wins", "Losses", "Ties") - copy formatting
# (bold, centered, bordered) from the existing header style (AC1) first,
# then set the text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record columns for every data row (2-53): Wins=71, Losses=91, Ties=0
$ws.Range("AD2:AD53").Value = 71
$ws.Range("AE2:AE53").Value = 91
$ws.Range("AF2:AF53").Value = 0
